# Daily attendance processing - 2025-12-01 14:56:43
# Normalizes the "Recorded By" column (G): whenever the recorder list
# begins with the literal token "System", that token is moved from the
# front of the comma-separated list to the back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1 -and $parts[0] -ceq "System") {
            $rest = $parts[1..($parts.Length - 1)]
            $newParts = $rest + @("System")
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value2 = $newVal
        }
    }
}
